$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1: for each of these row pairs, the two matches had been recorded with
# the team names (and all related odds/result columns) swapped. Fix it by
# exchanging columns F,G,H,I,J,L,M,N,P,Q,R,T,U,V between the two rows.
# (Columns A,B,C,D,E,K,O,S stay put - they are tied to the row position.)
# ---------------------------------------------------------------------------

function Swap-RowColumns($row1, $row2, $cols) {
    foreach ($col in $cols) {
        $v1 = $ws.Cells.Item($row1, $col).Value2
        $v2 = $ws.Cells.Item($row2, $col).Value2
        $ws.Cells.Item($row1, $col).Value = $v2
        $ws.Cells.Item($row2, $col).Value = $v1
    }
}

$swapCols = @(6, 7, 8, 9, 10, 12, 13, 14, 16, 17, 18, 20, 21, 22)

Swap-RowColumns 5 6 $swapCols
Swap-RowColumns 60 61 $swapCols
Swap-RowColumns 74 75 $swapCols
Swap-RowColumns 134 135 $swapCols

# ---------------------------------------------------------------------------
# Part 2: append 6 new match rows (138-143) at the bottom of the sheet,
# copying the formatting of the last existing row (137) first so the new
# rows pick up the same cell styles (bold/bordered index column, date format
# on the match-date column) used throughout the table.
# ---------------------------------------------------------------------------

$newRows = @(
    @{ A=137; B="portugal"; C="liga-portugal-2"; D="2023-2024"; E=45297.5;          F="Oliveirense";     G=1; H="Maritimo";      I=1; J=4.91; K="30/12/2024 15:12"; L=4.67; M="06/01/2024 11:33"; N=3.91; O="30/12/2024 15:12"; P=3.9;  Q="06/01/2024 11:49"; R=1.65; S="30/12/2024 15:12"; T=1.77; U="06/01/2024 11:49"; V="https://www.betexplorer.com/football/portugal/liga-portugal-2/oliveirense-maritimo/UFMkEDti/" },
    @{ A=138; B="portugal"; C="liga-portugal-2"; D="2023-2024"; E=45297.625;        F="Penafiel";        G=2; H="Vilaverdense";   I=1; J=1.85; K="30/12/2024 15:12"; L=1.69; M="06/01/2024 14:51"; N=3.44; O="30/12/2024 15:12"; P=3.72; Q="06/01/2024 14:50"; R=4.25; S="30/12/2024 15:12"; T=5.63; U="06/01/2024 14:54"; V="https://www.betexplorer.com/football/portugal/liga-portugal-2/penafiel-vilaverdense-fc/MBGZ7xKH/" },
    @{ A=139; B="portugal"; C="liga-portugal-2"; D="2023-2024"; E=45297.6875;       F="Benfica B";       G=2; H="Os Belenenses";  I=1; J=1.87; K="30/12/2024 18:13"; L=1.6;  M="06/01/2024 16:06"; N=3.61; O="30/12/2024 18:13"; P=4.07; Q="06/01/2024 16:27"; R=3.93; S="30/12/2024 18:13"; T=5.96; U="06/01/2024 16:08"; V="https://www.betexplorer.com/football/portugal/liga-portugal-2/benfica-cf-os-belenenses/t0T7kU3j/" },
    @{ A=140; B="portugal"; C="liga-portugal-2"; D="2023-2024"; E=45297.6875;       F="Santa Clara";     G=0; H="Mafra";          I=1; J=1.83; K="31/12/2024 12:12"; L=1.62; M="06/01/2024 16:21"; N=3.46; O="31/12/2024 12:12"; P=3.82; Q="06/01/2024 16:24"; R=4.31; S="31/12/2024 12:12"; T=6.39; U="06/01/2024 16:24"; V="https://www.betexplorer.com/football/portugal/liga-portugal-2/santa-clara-mafra/KnKIBKKh/" },
    @{ A=141; B="portugal"; C="liga-portugal-2"; D="2023-2024"; E=45297.79166666666;F="Academico Viseu"; G=1; H="Leiria";         I=0; J=2.33; K="30/12/2024 19:12"; L=2.51; M="06/01/2024 18:59"; N=3.27; O="30/12/2024 19:12"; P=3.18; Q="06/01/2024 18:55"; R=3.04; S="30/12/2024 19:12"; T=3.13; U="06/01/2024 18:59"; V="https://www.betexplorer.com/football/portugal/liga-portugal-2/academico-viseu-leiria/86LgDXec/" },
    @{ A=142; B="portugal"; C="liga-portugal-2"; D="2023-2024"; E=45297.79166666666;F="Nacional";        G=1; H="Tondela";         I=1; J=2.08; K="30/12/2024 15:12"; L=1.88; M="06/01/2024 16:57"; N=3.36; O="30/12/2024 15:12"; P=3.76; Q="06/01/2024 18:50"; R=3.5;  S="30/12/2024 15:12"; T=4.21; U="06/01/2024 18:50"; V="https://www.betexplorer.com/football/portugal/liga-portugal-2/nacional-tondela/29U3jAlp/" }
)

$firstNewRow = 138
$lastOldRow = 137
$lastNewRow = $firstNewRow + $newRows.Count - 1

# Copy the formatting (styles only) of the last existing row onto the
# range that will hold the new rows, so new cells inherit the right
# number formats / borders / bold index column.
$ws.Range("A$lastOldRow`:V$lastOldRow").Copy()
$ws.Range("A$firstNewRow`:V$lastNewRow").PasteSpecial(-4122)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $firstNewRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = $data.T
    $ws.Cells.Item($r, 21).Value = $data.U
    $ws.Cells.Item($r, 22).Value = $data.V
}

Write-Output "edit complete"
